$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8
$ws.Cells.Item($newRow, 1).Value = "{'activation': 'relu', 'alpha': 1e-05, 'beta_1': 0.9, 'hidden_layer_sizes': 10, 'learning_rate': 'constant', 'learning_rate_init': 0.1, 'max_iter': 500, 'momentum': 0.9, 'power_t': 0.5, 'random_state': 6, 'solver': 'lbfgs', 'tol': 1e-05}"
$ws.Cells.Item($newRow, 2).Value = 0.996
$ws.Cells.Item($newRow, 3).Value = 0.65
